# Update "想去人数" (attendance count) figures in the "展览" and "全部类型"
# sheets, reflecting refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first worksheet) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 5072
$ws1.Range("F5").Value = 7353
$ws1.Range("F10").Value = 66
$ws1.Range("F11").Value = 24
$ws1.Range("F12").Value = 4286
$ws1.Range("F13").Value = 1739
$ws1.Range("F15").Value = 98
$ws1.Range("F16").Value = 2891
$ws1.Range("F21").Value = 427
$ws1.Range("F23").Value = 292
$ws1.Range("F25").Value = 1682
$ws1.Range("F27").Value = 89
$ws1.Range("F28").Value = 1362
$ws1.Range("F35").Value = 105
$ws1.Range("F36").Value = 54
$ws1.Range("F37").Value = 2813
$ws1.Range("F38").Value = 698
$ws1.Range("F39").Value = 29

# --- Sheet "全部类型" (fourth worksheet) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 5072
$ws4.Range("F5").Value = 7353
$ws4.Range("F10").Value = 66
$ws4.Range("F11").Value = 24
$ws4.Range("F12").Value = 4286
$ws4.Range("F13").Value = 1739
$ws4.Range("F15").Value = 98
$ws4.Range("F16").Value = 2891
$ws4.Range("F21").Value = 427
$ws4.Range("F23").Value = 292
$ws4.Range("F25").Value = 1682
$ws4.Range("F27").Value = 89
$ws4.Range("F28").Value = 1362
$ws4.Range("F35").Value = 105
$ws4.Range("F36").Value = 54
$ws4.Range("F37").Value = 2813
$ws4.Range("F39").Value = 698
$ws4.Range("F40").Value = 29
